$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 272 (pushes existing rows 272..380 down to 273..381,
# and the style/format of the row above carries through for column D's date style).
$ws.Rows.Item(272).Insert()

# Fill in the new row 272 with the new data record.
$ws.Range("A272").Value = 3
$ws.Range("B272").Value = "Femacal de La Calera"
$ws.Range("C272").Value = "Coquimbo"
$ws.Range("D272").Value = 44755
$ws.Range("E272").Value = 5
$ws.Range("F272").Value = 100112043
$ws.Range("G272").Value = "Pepino ensalada"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 130
$ws.Range("K272").Value = 17000
$ws.Range("L272").Value = 18000
$ws.Range("M272").Value = 17538
$ws.Range("N272").Value = "$/caja 70 unidades"
$ws.Range("O272").Value = "Región de Arica y Parinacota"
$ws.Range("P272").Value = 251
$ws.Range("Q272").Value = 70
$ws.Range("R272").Value = "Hortaliza"
